# Scheduled runner update: refresh market-price derived columns (H-N) across Leve profit sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 3939.7273
$ws.Range("I6").Value = 269.4
$ws.Range("K6").Value = 808.1999999999999
$ws.Range("M6").Value = -696.1999999999999
# Row 98
$ws.Range("H98").Value = 8778.299999999999
$ws.Range("I98").Value = 8642.556
$ws.Range("K98").Value = 8642.556
$ws.Range("M98").Value = -7144.556
# Row 122
$ws.Range("H122").Value = 8778.299999999999
$ws.Range("I122").Value = 8642.556
$ws.Range("K122").Value = 25927.668
$ws.Range("M122").Value = -23477.668
# Row 129
$ws.Range("H129").Value = 2270.9285
$ws.Range("J129").Value = 2407.1538
$ws.Range("L129").Value = 7221.4614
$ws.Range("N129").Value = -17221.4614

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1265177.6
$ws.Range("I32").Value = 580232.9
$ws.Range("K32").Value = 580232.9
$ws.Range("M32").Value = -579945.9
# Row 74
$ws.Range("H74").Value = 2298.0789
$ws.Range("J74").Value = 3094.4443
$ws.Range("L74").Value = 3094.4443
$ws.Range("N74").Value = -4842.4443
# Row 77
$ws.Range("H77").Value = 2298.0789
$ws.Range("J77").Value = 3094.4443
$ws.Range("L77").Value = 15472.2215
$ws.Range("N77").Value = -24208.2215
# Row 88
$ws.Range("H88").Value = 3323.5833
$ws.Range("I88").Value = 2193.4
$ws.Range("J88").Value = 4130.857
$ws.Range("K88").Value = 2193.4
$ws.Range("L88").Value = 4130.857
$ws.Range("M88").Value = -1787.4
$ws.Range("N88").Value = -4942.857
# Row 91
$ws.Range("H91").Value = 3323.5833
$ws.Range("I91").Value = 2193.4
$ws.Range("J91").Value = 4130.857
$ws.Range("K91").Value = 2193.4
$ws.Range("L91").Value = 4130.857
$ws.Range("M91").Value = -789.4000000000001
$ws.Range("N91").Value = -6938.857
# Row 97
$ws.Range("H97").Value = 1051.2
$ws.Range("I97").Value = 1051.2
$ws.Range("K97").Value = 1051.2
$ws.Range("M97").Value = -555.2
# Row 132
$ws.Range("H132").Value = 1934.258
$ws.Range("I132").Value = 1160.65
$ws.Range("J132").Value = 3340.818
$ws.Range("K132").Value = 3481.95
$ws.Range("L132").Value = 10022.454
$ws.Range("M132").Value = -951.9500000000003
$ws.Range("N132").Value = -15082.454

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 100
$ws.Range("H100").Value = 24999.5
$ws.Range("J100").Value = 24999.5
$ws.Range("L100").Value = 24999.5
$ws.Range("N100").Value = -27163.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3294070.5
$ws.Range("I31").Value = 3139.1538
$ws.Range("J31").Value = 10424422
$ws.Range("K31").Value = 3139.1538
$ws.Range("L31").Value = 10424422
$ws.Range("M31").Value = -2844.1538
$ws.Range("N31").Value = -10425012
# Row 33
$ws.Range("H33").Value = 3755
$ws.Range("I33").Value = 3755
$ws.Range("K33").Value = 3755
$ws.Range("M33").Value = -3376
# Row 34
$ws.Range("H34").Value = 3294070.5
$ws.Range("I34").Value = 3139.1538
$ws.Range("J34").Value = 10424422
$ws.Range("K34").Value = 3139.1538
$ws.Range("L34").Value = 10424422
$ws.Range("M34").Value = -2937.1538
$ws.Range("N34").Value = -10424826
# Row 99
$ws.Range("H99").Value = 3052.889
$ws.Range("I99").Value = 2079.5
$ws.Range("J99").Value = 4999.6665
$ws.Range("K99").Value = 2079.5
$ws.Range("L99").Value = 4999.6665
$ws.Range("M99").Value = -581.5
$ws.Range("N99").Value = -7995.6665
# Row 126
$ws.Range("H126").Value = 3052.889
$ws.Range("I126").Value = 2079.5
$ws.Range("J126").Value = 4999.6665
$ws.Range("K126").Value = 6238.5
$ws.Range("L126").Value = 14998.9995
$ws.Range("M126").Value = -3768.5
$ws.Range("N126").Value = -19938.9995

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 168.8
$ws.Range("I7").Value = 168.8
$ws.Range("K7").Value = 506.4
$ws.Range("M7").Value = -394.4
# Row 40
$ws.Range("H40").Value = 237.6923
$ws.Range("I40").Value = 189.375
$ws.Range("J40").Value = 315
$ws.Range("K40").Value = 757.5
$ws.Range("L40").Value = 1260
$ws.Range("M40").Value = -688.5
$ws.Range("N40").Value = -1398
# Row 47
$ws.Range("H47").Value = 2177.818
$ws.Range("I47").Value = 864
$ws.Range("K47").Value = 2592
$ws.Range("M47").Value = -2161
# Row 68
$ws.Range("H68").Value = 12509038
$ws.Range("I68").Value = 700
$ws.Range("K68").Value = 2100
$ws.Range("M68").Value = -1289
# Row 71
$ws.Range("H71").Value = 12509038
$ws.Range("I71").Value = 700
$ws.Range("K71").Value = 6300
$ws.Range("M71").Value = -2244
# Row 80
$ws.Range("H80").Value = 19750
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
# Row 83
$ws.Range("H83").Value = 19750
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
# Row 92
$ws.Range("H92").Value = 639.8
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 9037.846
$ws.Range("J126").Value = 13258.75
$ws.Range("L126").Value = 39776.25
$ws.Range("N126").Value = -44716.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 8281.056
$ws.Range("I40").Value = 8597.308000000001
$ws.Range("K40").Value = 8597.308000000001
$ws.Range("M40").Value = -8461.308000000001
# Row 55
$ws.Range("H55").Value = 1109.091
$ws.Range("I55").Value = 897.6667
$ws.Range("J55").Value = 1362.8
$ws.Range("K55").Value = 897.6667
$ws.Range("L55").Value = 1362.8
$ws.Range("M55").Value = -724.6667
$ws.Range("N55").Value = -1708.8
# Row 93
$ws.Range("H93").Value = 1671.75
$ws.Range("I93").Value = 1267.8572
$ws.Range("J93").Value = 4499
$ws.Range("K93").Value = 1267.8572
$ws.Range("L93").Value = 4499
$ws.Range("M93").Value = -19.85719999999992
$ws.Range("N93").Value = -6995
# Row 122
$ws.Range("H122").Value = 9951.044
$ws.Range("I122").Value = 6810.091
$ws.Range("K122").Value = 20430.273
$ws.Range("M122").Value = -17980.273

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3643.4348
$ws.Range("I81").Value = 2912.625
$ws.Range("J81").Value = 5313.857
$ws.Range("K81").Value = 5825.25
$ws.Range("L81").Value = 10627.714
$ws.Range("M81").Value = -4764.25
$ws.Range("N81").Value = -12749.714
# Row 84
$ws.Range("H84").Value = 3643.4348
$ws.Range("I84").Value = 2912.625
$ws.Range("J84").Value = 5313.857
$ws.Range("K84").Value = 29126.25
$ws.Range("L84").Value = 53138.57
$ws.Range("M84").Value = -23822.25
$ws.Range("N84").Value = -63746.57
# Row 107
$ws.Range("H107").Value = 1063.2667
$ws.Range("J107").Value = 1826.2858
$ws.Range("L107").Value = 5478.857400000001
$ws.Range("N107").Value = -9318.857400000001
# Row 113
$ws.Range("H113").Value = 1506.1818
$ws.Range("I113").Value = 1046.5714
$ws.Range("K113").Value = 3139.7142
$ws.Range("M113").Value = -969.7142000000003
# Row 122
$ws.Range("H122").Value = 17860344
$ws.Range("I122").Value = 3402.25
$ws.Range("K122").Value = 10206.75
$ws.Range("M122").Value = -7756.75
# Row 126
$ws.Range("H126").Value = 1800.6
$ws.Range("I126").Value = 1499.5
$ws.Range("K126").Value = 4498.5
$ws.Range("M126").Value = -2028.5
# Row 136
$ws.Range("H136").Value = 4265.4814
$ws.Range("I136").Value = 3811.652
$ws.Range("K136").Value = 11434.956
$ws.Range("M136").Value = -8884.956

Write-Host "done"
